# Remove '[object object]' from bottom of empty attachments table
#
# The template text "{d.otherAttachments:ifEM():show(.noData)}" needs an
# ":elseShow()" clause appended right after ":show(.noData)" (and before the
# closing "}"), producing "{d.otherAttachments:ifEM():show(.noData):elseShow()}".
#
# The existing paragraph is made up of three runs:
#   [1] "{d."
#   [2] "otherAttachments:ifEM():show(.noData)"
#   [3] "}"
# We need to insert a brand-new fourth run containing ":elseShow()" between
# runs [2] and [3], while leaving runs [1]-[3] exactly as they were.
#
# Naively using Find/Replace or InsertAfter collapses/re-merges adjacent runs
# that share identical run formatting, which would wipe out the run
# boundaries. To stop that from happening we briefly give the runs on either
# side of the insertion point (the opening "{d." run and the closing "}" run)
# a distinguishing format (Bold), perform the insertion, and then restore
# their formatting back to normal - this keeps every run as its own distinct
# <w:r> element once the formatting is reverted.

$d = $word.ActiveDocument

$rngFull = $d.Content
$rngFull.Find.Execute("{d.otherAttachments:ifEM():show(.noData)}", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start = $rngFull.Start
$end = $rngFull.End

$newText = ":elseShow()"
$newLen = $newText.Length

# Temporarily mark the leading "{d." run (3 chars) as bold so it won't be
# re-merged with its neighbour once we start editing the paragraph.
$rngLead = $d.Range($start, $start + 3)
$rngLead.Bold = 1

# Temporarily mark the trailing "}" run (1 char) as bold too, for the same
# reason.
$rngTrail = $d.Range($end - 1, $end)
$rngTrail.Bold = 1

# Insert the new ":elseShow()" text right before the trailing "}" run (i.e.
# right after "...show(.noData)").
$rngInsertPoint = $d.Range($end - 1, $end - 1)
$rngInsertPoint.InsertAfter($newText)

# Mark the freshly inserted text bold as well so it stays its own run.
$rngNew = $d.Range($end - 1, $end - 1 + $newLen)
$rngNew.Bold = 1

# Now restore normal (non-bold) formatting on all three protected runs. Since
# each one is already its own distinct run, clearing the formatting keeps
# them separate instead of re-merging.
$rngLead2 = $d.Range($start, $start + 3)
$rngLead2.Bold = 0

$rngNew2 = $d.Range($end - 1, $end - 1 + $newLen)
$rngNew2.Bold = 0

$rngTrail2 = $d.Range($end - 1 + $newLen, $end + $newLen)
$rngTrail2.Bold = 0
